$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("F1:G82")
$rng.Replace("Breiten", "Breitensport", 1, 1, $false, $false, $false, $false)
$rng.Replace("Beginner", "Beginners", 1, 1, $false, $false, $false, $false)
$rng.Replace("Open", "Open Class", 1, 1, $false, $false, $false, $false)

$ws.Columns("F:G").AutoFit()
$ws.Range("G5").Select()
